$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 16.831
$ws.Range("B9").Value = 6.345999999999999
$ws.Range("E12").Value = 17.646
$ws.Range("B18").Value = 5.029999999999999
$ws.Range("B20").Value = 6.967000000000001
$ws.Range("E26").Value = 16.525
$ws.Range("B27").Value = 5.752000000000001
$ws.Range("E27").Value = 16.728
$ws.Range("E29").Value = 17
$ws.Range("E37").Value = 16.855
$ws.Range("E38").Value = 16.871
$ws.Range("E51").Value = 16.65
$ws.Range("E55").Value = 16.494
$ws.Range("B69").Value = 5.627
$ws.Range("E69").Value = 17.321
$ws.Range("E70").Value = 17.524
$ws.Range("B76").Value = 6.308
$ws.Range("B82").Value = 5.366000000000001
$ws.Range("E83").Value = 16.886
$ws.Range("E102").Value = 16.724
